# AssignData.xlsx - "Adding all the basic scenarios of Selenium"
#
# The hyperlinked e-mail/password test-data column (B1:B3) is converted back
# to plain data: the mailto: hyperlinks are removed, the cells lose the
# "Hyperlink" look (underline/blue), and the now-unused "Hyperlink" cell
# style is deleted from the workbook. The B2 entry is replaced with a plain
# numeric OTP-style value (123456) instead of the repeated "Queen@123" text,
# and the active selection moves from B3 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the three mailto: hyperlinks that lived on B1:B3.
$ws.Hyperlinks.Delete()

# Strip the inherited "Hyperlink" formatting (underline, blue font) from
# those cells so they fall back to the default/Normal look.
$ws.Range("B1:B3").ClearFormats()

# The "Hyperlink" cell style is no longer used anywhere - drop it.
$wb.Styles.Item("Hyperlink").Delete()

# B2 becomes a plain numeric value instead of the shared "Queen@123" text.
$ws.Range("B2").Value = 123456

# Move the active selection/cell from B3 to B2.
$ws.Range("B2").Select() | Out-Null
